# Auto-generated: apply cryptos list update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @("D2", "26.675.75"),
    @("E2", "  +4.33%  "),
    @("D3", "1.753.46"),
    @("E4", "  -0.13%  "),
    @("D5", "247.47"),
    @("E5", "  +3.55%  "),
    @("D6", "0.9982"),
    @("E6", "  -0.15%  "),
    @("D7", "0.4813"),
    @("E7", "  +0.02%  "),
    @("D8", "0.2719"),
    @("E8", "  +3.33%  "),
    @("D9", "0.06265"),
    @("E9", "  +1.26%  "),
    @("D10", "1.744.39"),
    @("E10", "  +4.50%  "),
    @("D11", "0.07122"),
    @("E11", "  +1.51%  "),
    @("D12", "15.97"),
    @("E12", "  +7.02%  "),
    @("D13", "0.6258"),
    @("E13", "  +5.62%  "),
    @("D14", "4.524"),
    @("E14", "  +3.01%  "),
    @("D15", "77.43"),
    @("E15", "  +2.79%  "),
    @("D16", "0.9982"),
    @("E16", "  -0.15%  "),
    @("D17", "26.667.20"),
    @("E17", "  +4.31%  "),
    @("D18", "0.9981"),
    @("E18", "  -0.13%  "),
    @("D19", "0.000006914"),
    @("E19", "  +1.90%  "),
    @("D20", "11.79"),
    @("D21", "1.966.74"),
    @("E21", "  +4.49%  "),
    @("D22", "4.663"),
    @("E22", "  +4.73%  "),
    @("D23", "8.895"),
    @("E23", "  +1.77%  "),
    @("D24", "5.365"),
    @("E24", "  +1.35%  "),
    @("D25", "136.60"),
    @("E25", "  -0.21%  "),
    @("D26", "15.52"),
    @("E26", "  +2.94%  "),
    @("D27", "1.844"),
    @("E27", "  +6.55%  "),
    @("D28", "1.416"),
    @("E28", "  +1.64%  "),
    @("D29", "108.06"),
    @("E29", "  +2.99%  "),
    @("D30", "4.051"),
    @("E30", "  +1.38%  "),
    @("D31", "3.782"),
    @("E31", "  +3.39%  "),
    @("E32", "  +1.20%  "),
    @("D33", "0.04610"),
    @("E33", "  +8.78%  "),
    @("D34", "2.610"),
    @("E34", "  -0.24%  "),
    @("B35", "ImmutableX"),
    @("C35", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"),
    @("D35", "0.6420"),
    @("E35", "  +5.13%  "),
    @("B36", "ARBITRUM"),
    @("C36", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"),
    @("D36", "1.007"),
    @("E36", "  +5.81%  "),
    @("D37", "0.9488"),
    @("E37", "  +10.39%  "),
    @("D38", "114.65"),
    @("E38", "  +19.18%  "),
    @("D39", "2.508"),
    @("E39", "  -3.40%  "),
    @("D40", "2.006"),
    @("E40", "  +7.50%  "),
    @("E41", "  +0.34%  "),
    @("D42", "0.01519"),
    @("E42", "  +2.98%  "),
    @("D43", "5.760"),
    @("E43", "  +18.76%  "),
    @("D44", "0.3937"),
    @("E44", "  +4.18%  "),
    @("D45", "6.770"),
    @("E45", "  +8.73%  "),
    @("D46", "0.1208"),
    @("E46", "  +8.01%  "),
    @("E47", "  +1.63%  "),
    @("D48", "8.001"),
    @("E48", "  +8.63%  "),
    @("D49", "31.05"),
    @("E49", "  +3.90%  "),
    @("D50", "0.3472"),
    @("E50", "  +3.71%  "),
    @("D51", "52.00"),
    @("E51", "  +3.70%  ")
)

foreach ($u in $updates) {
    $addr = $u[0]
    $val = $u[1]
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Write-Host "Applied $($updates.Count) cell updates"
